$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" sheet by duplicating the "2022-Q3" sheet
#    (so it inherits identical column layout / header styling), placing
#    the copy immediately before "2022-Q3" (i.e. right after "总计").
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Drop the old Q3 data rows (3-9) that came along with the copy; only one
# fund is reported for 2022-Q4.
$q4.Range("A3:A9").EntireRow.Delete()

# Overwrite row 2 with the 2022-Q4 fund holding.
$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = "'007592"
$q4.Cells.Item(2,3).Value = "华夏价值精选混合"
$q4.Cells.Item(2,4).Value = "'4.06"
$q4.Cells.Item(2,5).Value = "'93.21"
$q4.Cells.Item(2,6).Value = "'3.05"
$q4.Cells.Item(2,7).Value = "'0.1238"
$q4.Cells.Item(2,8).Value = 9

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: shift the existing Q3/Q2 rows down
#    by one and insert the new Q4 summary row above them.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Shift old row 3 (2022-Q2) down to row 4.
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = $total.Cells.Item(3,2).Value2
$total.Cells.Item(4,3).Value = $total.Cells.Item(3,3).Value2
$total.Cells.Item(4,4).Value = $total.Cells.Item(3,4).Value2
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(4,1).PasteSpecial(-4122)

# Shift old row 2 (2022-Q3) down to row 3.
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = $total.Cells.Item(2,2).Value2
$total.Cells.Item(3,3).Value = $total.Cells.Item(2,3).Value2
$total.Cells.Item(3,4).Value = $total.Cells.Item(2,4).Value2

# Write the new 2022-Q4 summary row into row 2.
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 1
$total.Cells.Item(2,4).Value = 0.12

# Restore the originally-active sheet/selection (the new-sheet copy steals
# tab focus), matching the source workbook's selection state.
$wb.Worksheets.Item("2022-Q2").Activate()

Write-Host "2022-Q4 sheet added"
